# CategoryWeight_CalcTool.xlsx edit — convert the category list (usa /
# australia / india) into the seven days of the week, update the backing
# weight numbers (row 4-10) accordingly, extend the SUM/aggregate ranges
# down to row 19, add the stray I19 helper value, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Category names (column A) -------------------------------------------
$ws.Range("A4").Value  = "mon"
$ws.Range("A5").Value  = "tues"
$ws.Range("A6").Value  = "wed"
$ws.Range("A7").Value  = "thurs"
$ws.Range("A8").Value  = "fri"
$ws.Range("A9").Value  = "sat"
$ws.Range("A10").Value = "sun"

# --- Weights (column B) ---------------------------------------------------
$ws.Range("B4").Value     = 1
$ws.Range("B5").Formula   = "=B4*1.25"
$ws.Range("B6").Formula   = "=B5*1.25"
$ws.Range("B7").Value     = 1
$ws.Range("B8").Value     = 1
$ws.Range("B9").Value     = 1.25
$ws.Range("B10").Value    = 1.25

# --- Aggregate formulas that now span the full category block ------------
$ws.Range("B22").Formula = "=SUM(B4:B19)"
$ws.Range("H19").Formula = "=D9"

# --- Stray helper value introduced alongside the new rows -----------------
$ws.Range("I19").Value = 18

# --- Selection / view ------------------------------------------------------
$ws.Range("B7:B10").Select()

$wb.Application.Calculate()
